# Upgrade Conda env for Geopandas to Conda Forge - reshape venues sheet:
# insert two new "Unnamed: 0.1" / "Unnamed: 0" index columns before venue_id,
# shifting the original venue_id..geometry columns two slots to the right,
# and normalize empty-geometry text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find how many rows are used (38 in the source data).
$lastRow = $ws.UsedRange.Rows.Count

# Insert two new columns at C and D; this pushes the existing
# venue..geometry columns (C..L) to E..N. Column B does not shift
# (it sits to the left of the insertion point), so its old header
# ("venue_id") must be copied over to D explicitly afterwards.
$oldB1 = $ws.Range("B1").Value2
$ws.Range("C1:D1").EntireColumn.Insert()

# New header labels: B/C get the two new synthetic index headers,
# and D gets back the original "venue_id" header that used to live
# in B1 before the insert.
$ws.Range("B1").Value = "Unnamed: 0.1"
$ws.Range("C1").Value = "Unnamed: 0"
$ws.Range("D1").Value = $oldB1

# For the trailing summary rows, column B (venue_id) was blank in the
# source; mirror column A's row index into B before propagating it to
# the two newly inserted columns, matching the rest of the data rows.
for ($r = 2; $r -le $lastRow; $r++) {
    $bVal = $ws.Cells.Item($r, 2).Value2
    if ($bVal -eq $null -or $bVal -eq "") {
        $ws.Cells.Item($r, 2).Value = $ws.Cells.Item($r, 1).Value2
    }
}

# Populate the two new columns (C, D) with the same values as column B
# for every data row, mirroring a pandas reset_index producing
# "Unnamed: 0.1"/"Unnamed: 0" duplicate index columns.
for ($r = 2; $r -le $lastRow; $r++) {
    $bVal = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 3).Value = $bVal
    $ws.Cells.Item($r, 4).Value = $bVal
}

# The empty-geometry rows previously rendered as "POINT (nan nan)" now
# render as "POINT EMPTY" (now in column N after the column insert).
for ($r = 2; $r -le $lastRow; $r++) {
    $g = $ws.Cells.Item($r, 14).Value2
    if ($g -eq "POINT (nan nan)") {
        $ws.Cells.Item($r, 14).Value = "POINT EMPTY"
    }
}
